$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "Upload backend code in S3 bucket." paragraph
#   becomes: Upload backend code in S3 bucket name "backend-userservices".
#   with the bold run: backend-userservices".
# ---------------------------------------------------------------------------

$target = "Upload backend code in S3 bucket."
$found = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $ptext = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($ptext -eq $target) {
        $found = $p
        break
    }
}

if ($found -ne $null) {
    $r = $found.Range
    $prefix = "Upload backend code in S3 bucket"

    # Replace whole paragraph text (minus the trailing paragraph mark) with
    # the new un-bolded lead-in text.
    $plainRange = $d.Range($r.Start, $r.End - 1)
    $plainRange.Text = $prefix

    $pos = $r.Start + $prefix.Length
    $ip = $d.Range($pos, $pos)
    $ip.InsertAfter(" name " + [char]0x201C)

    $pos = $pos + (" name " + [char]0x201C).Length
    $ip = $d.Range($pos, $pos)
    $ip.InsertAfter("backend-userservices")
    $ip.Font.Bold = 1

    $pos = $pos + ("backend-userservices").Length
    $ip = $d.Range($pos, $pos)
    $ip.InsertAfter([char]0x201D)
    $ip.Font.Bold = 1

    $pos = $pos + 1
    $ip = $d.Range($pos, $pos)
    $ip.InsertAfter(".")
    $ip.Font.Bold = 1
}

# ---------------------------------------------------------------------------
# Edit 2: "Upload build in frontend bucket." paragraph
#   becomes: Upload build in frontend S3 bucket.
# ---------------------------------------------------------------------------

$ok = $d.Content.Find.Execute("in frontend bucket.", $true, $false, $false, `
    $false, $false, $true, 1, $false, "in frontend S3 bucket.", 2)
